$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 63, shifting the existing rows 63-65 down to 64-66.
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new weekly price record.
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 44448
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100102
$ws.Range("H63").Value = "Cítricos"
$ws.Range("I63").Value = 100102004
$ws.Range("J63").Value = "Mandarina"
$ws.Range("K63").Value = "Murcott"
$ws.Range("L63").Value = "Segunda"
$ws.Range("M63").Value = 250
$ws.Range("N63").Value = 12000
$ws.Range("O63").Value = 13000
$ws.Range("P63").Value = 12500
$ws.Range("Q63").Value = "$/caja 20 kilos"
$ws.Range("R63").Value = "Región de Coquimbo"
$ws.Range("S63").Value = 625
$ws.Range("T63").Value = 20
